$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: new period headers (match existing header style from N1)
$ws.Range("O1").Value = "31/03/2024"
$ws.Range("P1").Value = "30/06/2024"
$ws.Range("O1:P1").Font.Bold = $true
$ws.Range("O1:P1").HorizontalAlignment = -4108
$ws.Range("O1:P1").VerticalAlignment = -4160
$ws.Range("O1:P1").Borders.LineStyle = 1

$ws.Range("O2").Value = 7442510.848
$ws.Range("P2").Value = 8530219.007999999
$ws.Range("O3").Value = 4849333.248
$ws.Range("P3").Value = 5878986.24
$ws.Range("O4").Value = 935238.0159999999
$ws.Range("P4").Value = 1448705.024
$ws.Range("O5").Value = 146480
$ws.Range("P5").Value = 185148.992
$ws.Range("O6").Value = 1596228.992
$ws.Range("P6").Value = 798339.968
$ws.Range("O7").Value = 1857696
$ws.Range("P7").Value = 3086614.016
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("O9").Value = 220938
$ws.Range("P9").Value = 213212
$ws.Range("O10").Value = 11495
$ws.Range("P10").Value = 5447
$ws.Range("O11").Value = 81257
$ws.Range("P11").Value = 141519.008
$ws.Range("O12").Value = 462947.008
$ws.Range("P12").Value = 416403.008
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("O16").Value = 7816
$ws.Range("P16").Value = 6000
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("O19").Value = 257468.992
$ws.Range("P19").Value = 192968.992
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("O22").Value = 1792
$ws.Range("P22").Value = 5183
$ws.Range("O23").Value = 2100045.952
$ws.Range("P23").Value = 2188395.008
$ws.Range("O24").Value = 28393
$ws.Range("P24").Value = 41252
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("O26").Value = 7442510.848
$ws.Range("P26").Value = 8530219.007999999
$ws.Range("O27").Value = 3164288
$ws.Range("P27").Value = 3569024
$ws.Range("O28").Value = 38205
$ws.Range("P28").Value = 47882
$ws.Range("O29").Value = 2323632.128
$ws.Range("P29").Value = 2430176
$ws.Range("O30").Value = 19285
$ws.Range("P30").Value = 17842
$ws.Range("O31").Value = 608856
$ws.Range("P31").Value = 649752
$ws.Range("O32").Value = 0
$ws.Range("P32").Value = 0
$ws.Range("O33").Value = 0
$ws.Range("P33").Value = 0
$ws.Range("O34").Value = 174310
$ws.Range("P34").Value = 423372
$ws.Range("O35").Value = 0
$ws.Range("P35").Value = 0
$ws.Range("O36").Value = 0
$ws.Range("P36").Value = 0
$ws.Range("O37").Value = 838478.0159999999
$ws.Range("P37").Value = 1360870.016
$ws.Range("O38").Value = 800628.992
$ws.Range("P38").Value = 1322889.984
$ws.Range("O39").Value = 0
$ws.Range("P39").Value = 0
$ws.Range("O40").Value = 24126
$ws.Range("P40").Value = 23559
$ws.Range("O41").Value = 0
$ws.Range("P41").Value = 0
$ws.Range("O42").Value = 0
$ws.Range("P42").Value = 0
$ws.Range("O43").Value = 13723
$ws.Range("P43").Value = 14421
$ws.Range("O44").Value = 0
$ws.Range("P44").Value = 0
$ws.Range("O45").Value = 0
$ws.Range("P45").Value = 0
$ws.Range("O46").Value = 4523
$ws.Range("P46").Value = 5691
$ws.Range("O47").Value = 3435222.024
$ws.Range("P47").Value = 3594634.12
$ws.Range("O48").Value = 1518662.016
$ws.Range("P48").Value = 1518662.016
$ws.Range("O49").Value = 33847
$ws.Range("P49").Value = 35252
$ws.Range("O50").Value = 0
$ws.Range("P50").Value = 0
$ws.Range("O51").Value = 1880206.976
$ws.Range("P51").Value = 2034727.04
$ws.Range("O52").Value = 0
$ws.Range("P52").Value = 0
$ws.Range("O53").Value = 1768
$ws.Range("P53").Value = 1532
$ws.Range("O54").Value = 738
$ws.Range("P54").Value = 4461
$ws.Range("O55").Value = 0
$ws.Range("P55").Value = 0
$ws.Range("O56").Value = 0
$ws.Range("P56").Value = 0
# Row 57: blank separator row - force empty cells to exist at O57/P57
$ws.Range("O57:P57").NumberFormat = "General"
$ws.Range("O57").Value = ""
$ws.Range("P57").Value = ""
# Row 58: blank separator row - force empty cells to exist at O58/P58
$ws.Range("O58:P58").NumberFormat = "General"
$ws.Range("O58").Value = ""
$ws.Range("P58").Value = ""
$ws.Range("O59").Value = 2679222.016
$ws.Range("P59").Value = 2796468.992
$ws.Range("O60").Value = -2178440.96
$ws.Range("P60").Value = -2080269.056
$ws.Range("O61").Value = 500780.992
$ws.Range("P61").Value = 716200
$ws.Range("O62").Value = -303017.984
$ws.Range("P62").Value = -298720
$ws.Range("O63").Value = -18957
$ws.Range("P63").Value = -21974
$ws.Range("O64").Value = -1172
$ws.Range("P64").Value = -13093
$ws.Range("O65").Value = 8270
$ws.Range("P65").Value = -3731
$ws.Range("O66").Value = 0
$ws.Range("P66").Value = 0
$ws.Range("O67").Value = 0
$ws.Range("P67").Value = -68
$ws.Range("O68").Value = -9064
$ws.Range("P68").Value = -162104.992
$ws.Range("O69").Value = 90811
$ws.Range("P69").Value = 157088.992
$ws.Range("O70").Value = -99875
$ws.Range("P70").Value = -319193.984
# Row 71: blank separator row - force empty cells to exist at O71/P71
$ws.Range("O71:P71").NumberFormat = "General"
$ws.Range("O71").Value = ""
$ws.Range("P71").Value = ""
# Row 72: blank separator row - force empty cells to exist at O72/P72
$ws.Range("O72:P72").NumberFormat = "General"
$ws.Range("O72").Value = ""
$ws.Range("P72").Value = ""
# Row 73: blank separator row - force empty cells to exist at O73/P73
$ws.Range("O73:P73").NumberFormat = "General"
$ws.Range("O73").Value = ""
$ws.Range("P73").Value = ""
$ws.Range("O74").Value = 176840
$ws.Range("P74").Value = 216508.992
$ws.Range("O75").Value = -4152
$ws.Range("P75").Value = -6420
$ws.Range("O76").Value = -16250
$ws.Range("P76").Value = -64445
# Row 77: blank separator row - force empty cells to exist at O77/P77
$ws.Range("O77:P77").NumberFormat = "General"
$ws.Range("O77").Value = ""
$ws.Range("P77").Value = ""
# Row 78: blank separator row - force empty cells to exist at O78/P78
$ws.Range("O78:P78").NumberFormat = "General"
$ws.Range("O78").Value = ""
$ws.Range("P78").Value = ""
$ws.Range("O79").Value = 986
$ws.Range("P79").Value = 634
$ws.Range("O80").Value = 157424
$ws.Range("P80").Value = 146278

Write-Host "Applied TTEN3 Q1/Q2 2024 columns"
